$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column (H), copying the formatting used by the
# existing header cells (e.g. G1: bold, centered, bordered) so the new
# header matches the rest of row 1.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column values for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
